# Applies cyclic rotation of match-data rows that were reordered in the source
# (each group of rows sharing the same Date got its non-index columns rotated).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Group: rows [38, 39]
$ws.Cells.Item(38, 2).Value = 6155944
$ws.Cells.Item(38, 5).Value = 'Cobreloa'
$ws.Cells.Item(38, 6).Value = 'Puerto Montt'
$ws.Cells.Item(38, 8).Value = 0
$ws.Cells.Item(38, 9).Value = 1
$ws.Cells.Item(38, 10).Value = 0
$ws.Cells.Item(38, 11).Value = 'H'
$ws.Cells.Item(38, 12).Value = 1.4
$ws.Cells.Item(38, 13).Value = 3.75
$ws.Cells.Item(38, 14).Value = 7.5
$ws.Cells.Item(38, 15).Value = 1.45
$ws.Cells.Item(38, 16).Value = 3.6
$ws.Cells.Item(38, 17).Value = 6.5
$ws.Cells.Item(38, 18).Value = -1.25
$ws.Cells.Item(38, 19).Value = 2.05
$ws.Cells.Item(38, 20).Value = 1.75
$ws.Cells.Item(38, 22).Value = 1.9
$ws.Cells.Item(38, 23).Value = 1.9
$ws.Cells.Item(38, 24).Value = 0.45
$ws.Cells.Item(38, 26).Value = -1
$ws.Cells.Item(38, 27).Value = -0.5
$ws.Cells.Item(38, 28).Value = 0.375
$ws.Cells.Item(38, 29).Value = -1
$ws.Cells.Item(38, 30).Value = 0.8999999999999999
$ws.Cells.Item(39, 2).Value = 6158713
$ws.Cells.Item(39, 5).Value = 'Santiago Morning'
$ws.Cells.Item(39, 6).Value = 'CD Antofagasta'
$ws.Cells.Item(39, 8).Value = 2
$ws.Cells.Item(39, 9).Value = 0
$ws.Cells.Item(39, 10).Value = 1
$ws.Cells.Item(39, 11).Value = 'A'
$ws.Cells.Item(39, 12).Value = 3
$ws.Cells.Item(39, 13).Value = 3.25
$ws.Cells.Item(39, 14).Value = 2.1
$ws.Cells.Item(39, 15).Value = 2.625
$ws.Cells.Item(39, 16).Value = 3.2
$ws.Cells.Item(39, 17).Value = 2.3
$ws.Cells.Item(39, 18).Value = 0
$ws.Cells.Item(39, 19).Value = 2.025
$ws.Cells.Item(39, 20).Value = 1.775
$ws.Cells.Item(39, 22).Value = 1.75
$ws.Cells.Item(39, 23).Value = 1.95
$ws.Cells.Item(39, 24).Value = -1
$ws.Cells.Item(39, 26).Value = 1.3
$ws.Cells.Item(39, 27).Value = -1
$ws.Cells.Item(39, 28).Value = 0.7749999999999999
$ws.Cells.Item(39, 29).Value = 0.75
$ws.Cells.Item(39, 30).Value = -1

# Group: rows [43, 44]
$ws.Cells.Item(43, 2).Value = 6393620
$ws.Cells.Item(43, 5).Value = 'Union San Felipe'
$ws.Cells.Item(43, 6).Value = 'Santiago Morning'
$ws.Cells.Item(43, 8).Value = 1
$ws.Cells.Item(43, 9).Value = 1
$ws.Cells.Item(43, 11).Value = 'H'
$ws.Cells.Item(43, 12).Value = 1.95
$ws.Cells.Item(43, 13).Value = 3
$ws.Cells.Item(43, 14).Value = 3.75
$ws.Cells.Item(43, 15).Value = 1.65
$ws.Cells.Item(43, 17).Value = 5
$ws.Cells.Item(43, 18).Value = -0.75
$ws.Cells.Item(43, 19).Value = 1.9
$ws.Cells.Item(43, 20).Value = 1.9
$ws.Cells.Item(43, 21).Value = 2.5
$ws.Cells.Item(43, 22).Value = 1.9
$ws.Cells.Item(43, 23).Value = 1.9
$ws.Cells.Item(43, 24).Value = 0.6499999999999999
$ws.Cells.Item(43, 25).Value = -1
$ws.Cells.Item(43, 27).Value = 0.45
$ws.Cells.Item(43, 28).Value = -0.5
$ws.Cells.Item(43, 29).Value = 0.8999999999999999
$ws.Cells.Item(44, 2).Value = 6155945
$ws.Cells.Item(44, 5).Value = 'Puerto Montt'
$ws.Cells.Item(44, 6).Value = 'Deportes Iquique'
$ws.Cells.Item(44, 8).Value = 2
$ws.Cells.Item(44, 9).Value = 2
$ws.Cells.Item(44, 11).Value = 'D'
$ws.Cells.Item(44, 12).Value = 2.625
$ws.Cells.Item(44, 13).Value = 3.2
$ws.Cells.Item(44, 14).Value = 2.375
$ws.Cells.Item(44, 15).Value = 2.6
$ws.Cells.Item(44, 17).Value = 2.375
$ws.Cells.Item(44, 18).Value = 0
$ws.Cells.Item(44, 19).Value = 2
$ws.Cells.Item(44, 20).Value = 1.8
$ws.Cells.Item(44, 21).Value = 2.25
$ws.Cells.Item(44, 22).Value = 1.75
$ws.Cells.Item(44, 23).Value = 1.95
$ws.Cells.Item(44, 24).Value = -1
$ws.Cells.Item(44, 25).Value = 2.25
$ws.Cells.Item(44, 27).Value = 0
$ws.Cells.Item(44, 28).Value = 0
$ws.Cells.Item(44, 29).Value = 0.75

# Group: rows [61, 62]
$ws.Cells.Item(61, 2).Value = 6156906
$ws.Cells.Item(61, 5).Value = 'Barnechea'
$ws.Cells.Item(61, 6).Value = 'Deportes Temuco'
$ws.Cells.Item(61, 7).Value = 1
$ws.Cells.Item(61, 8).Value = 2
$ws.Cells.Item(61, 10).Value = 1
$ws.Cells.Item(61, 11).Value = 'A'
$ws.Cells.Item(61, 12).Value = 2.2
$ws.Cells.Item(61, 13).Value = 3
$ws.Cells.Item(61, 14).Value = 3.1
$ws.Cells.Item(61, 15).Value = 1.95
$ws.Cells.Item(61, 16).Value = 3.1
$ws.Cells.Item(61, 17).Value = 3.75
$ws.Cells.Item(61, 18).Value = -0.25
$ws.Cells.Item(61, 19).Value = 1.725
$ws.Cells.Item(61, 20).Value = 1.975
$ws.Cells.Item(61, 21).Value = 2.25
$ws.Cells.Item(61, 22).Value = 1.8
$ws.Cells.Item(61, 23).Value = 2
$ws.Cells.Item(61, 24).Value = -1
$ws.Cells.Item(61, 26).Value = 2.75
$ws.Cells.Item(61, 27).Value = -1
$ws.Cells.Item(61, 28).Value = 0.9750000000000001
$ws.Cells.Item(61, 29).Value = 0.8
$ws.Cells.Item(61, 30).Value = -1
$ws.Cells.Item(62, 2).Value = 6156905
$ws.Cells.Item(62, 5).Value = 'San Marcos De Arica'
$ws.Cells.Item(62, 6).Value = 'Rangers de Talca'
$ws.Cells.Item(62, 7).Value = 2
$ws.Cells.Item(62, 8).Value = 1
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 'H'
$ws.Cells.Item(62, 12).Value = 2
$ws.Cells.Item(62, 13).Value = 3.2
$ws.Cells.Item(62, 14).Value = 3.3
$ws.Cells.Item(62, 15).Value = 1.8
$ws.Cells.Item(62, 16).Value = 3.3
$ws.Cells.Item(62, 17).Value = 3.8
$ws.Cells.Item(62, 18).Value = -0.5
$ws.Cells.Item(62, 19).Value = 1.875
$ws.Cells.Item(62, 20).Value = 1.925
$ws.Cells.Item(62, 21).Value = 2.75
$ws.Cells.Item(62, 22).Value = 1.975
$ws.Cells.Item(62, 23).Value = 1.825
$ws.Cells.Item(62, 24).Value = 0.8
$ws.Cells.Item(62, 26).Value = -1
$ws.Cells.Item(62, 27).Value = 0.875
$ws.Cells.Item(62, 28).Value = -1
$ws.Cells.Item(62, 29).Value = 0.4875
$ws.Cells.Item(62, 30).Value = -0.5

# Group: rows [73, 74]
$ws.Cells.Item(73, 2).Value = 6156180
$ws.Cells.Item(73, 5).Value = 'CD Antofagasta'
$ws.Cells.Item(73, 6).Value = 'Barnechea'
$ws.Cells.Item(73, 7).Value = 2
$ws.Cells.Item(73, 8).Value = 0
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 11).Value = 'H'
$ws.Cells.Item(73, 12).Value = 1.615
$ws.Cells.Item(73, 13).Value = 3.6
$ws.Cells.Item(73, 14).Value = 4.5
$ws.Cells.Item(73, 15).Value = 1.666
$ws.Cells.Item(73, 16).Value = 3.75
$ws.Cells.Item(73, 17).Value = 4
$ws.Cells.Item(73, 18).Value = -0.75
$ws.Cells.Item(73, 19).Value = 1.875
$ws.Cells.Item(73, 20).Value = 1.925
$ws.Cells.Item(73, 21).Value = 2.5
$ws.Cells.Item(73, 22).Value = 1.925
$ws.Cells.Item(73, 23).Value = 1.875
$ws.Cells.Item(73, 24).Value = 0.6659999999999999
$ws.Cells.Item(73, 26).Value = -1
$ws.Cells.Item(73, 27).Value = 0.875
$ws.Cells.Item(73, 28).Value = -1
$ws.Cells.Item(73, 29).Value = -1
$ws.Cells.Item(73, 30).Value = 0.875
$ws.Cells.Item(74, 2).Value = 6155954
$ws.Cells.Item(74, 5).Value = 'La Serena'
$ws.Cells.Item(74, 6).Value = 'Puerto Montt'
$ws.Cells.Item(74, 7).Value = 1
$ws.Cells.Item(74, 8).Value = 2
$ws.Cells.Item(74, 9).Value = 1
$ws.Cells.Item(74, 11).Value = 'A'
$ws.Cells.Item(74, 12).Value = 1.75
$ws.Cells.Item(74, 13).Value = 3.3
$ws.Cells.Item(74, 14).Value = 4.2
$ws.Cells.Item(74, 15).Value = 1.909
$ws.Cells.Item(74, 16).Value = 3.2
$ws.Cells.Item(74, 17).Value = 3.75
$ws.Cells.Item(74, 18).Value = -0.5
$ws.Cells.Item(74, 19).Value = 1.95
$ws.Cells.Item(74, 20).Value = 1.85
$ws.Cells.Item(74, 21).Value = 2.25
$ws.Cells.Item(74, 22).Value = 2.025
$ws.Cells.Item(74, 23).Value = 1.775
$ws.Cells.Item(74, 24).Value = -1
$ws.Cells.Item(74, 26).Value = 2.75
$ws.Cells.Item(74, 27).Value = -1
$ws.Cells.Item(74, 28).Value = 0.8500000000000001
$ws.Cells.Item(74, 29).Value = 1.025
$ws.Cells.Item(74, 30).Value = -1

# Group: rows [101, 102, 103]
$ws.Cells.Item(101, 2).Value = 6156929
$ws.Cells.Item(101, 5).Value = 'Cobreloa'
$ws.Cells.Item(101, 6).Value = 'Union San Felipe'
$ws.Cells.Item(101, 8).Value = 0
$ws.Cells.Item(101, 10).Value = 0
$ws.Cells.Item(101, 11).Value = 'H'
$ws.Cells.Item(101, 12).Value = 2.3
$ws.Cells.Item(101, 13).Value = 2.9
$ws.Cells.Item(101, 14).Value = 3
$ws.Cells.Item(101, 15).Value = 1.8
$ws.Cells.Item(101, 17).Value = 4
$ws.Cells.Item(101, 18).Value = -0.5
$ws.Cells.Item(101, 19).Value = 1.875
$ws.Cells.Item(101, 20).Value = 1.925
$ws.Cells.Item(101, 22).Value = 1.9
$ws.Cells.Item(101, 23).Value = 1.9
$ws.Cells.Item(101, 24).Value = 0.8
$ws.Cells.Item(101, 25).Value = -1
$ws.Cells.Item(101, 27).Value = 0.875
$ws.Cells.Item(101, 28).Value = -1
$ws.Cells.Item(101, 30).Value = 0.8999999999999999
$ws.Cells.Item(102, 2).Value = 6156183
$ws.Cells.Item(102, 5).Value = 'Santiago Wanderers'
$ws.Cells.Item(102, 6).Value = 'CD Antofagasta'
$ws.Cells.Item(102, 7).Value = 3
$ws.Cells.Item(102, 9).Value = 1
$ws.Cells.Item(102, 12).Value = 2.2
$ws.Cells.Item(102, 13).Value = 3.2
$ws.Cells.Item(102, 15).Value = 2.2
$ws.Cells.Item(102, 16).Value = 3.2
$ws.Cells.Item(102, 17).Value = 3
$ws.Cells.Item(102, 18).Value = -0.25
$ws.Cells.Item(102, 19).Value = 1.95
$ws.Cells.Item(102, 20).Value = 1.85
$ws.Cells.Item(102, 21).Value = 2.25
$ws.Cells.Item(102, 22).Value = 1.8
$ws.Cells.Item(102, 23).Value = 2
$ws.Cells.Item(102, 24).Value = 1.2
$ws.Cells.Item(102, 27).Value = 0.95
$ws.Cells.Item(102, 29).Value = 0.8
$ws.Cells.Item(102, 30).Value = -1
$ws.Cells.Item(103, 2).Value = 6155958
$ws.Cells.Item(103, 5).Value = 'Puerto Montt'
$ws.Cells.Item(103, 6).Value = 'San Marcos De Arica'
$ws.Cells.Item(103, 7).Value = 1
$ws.Cells.Item(103, 8).Value = 1
$ws.Cells.Item(103, 9).Value = 0
$ws.Cells.Item(103, 10).Value = 1
$ws.Cells.Item(103, 11).Value = 'D'
$ws.Cells.Item(103, 12).Value = 3
$ws.Cells.Item(103, 14).Value = 2.2
$ws.Cells.Item(103, 15).Value = 2.05
$ws.Cells.Item(103, 16).Value = 3.3
$ws.Cells.Item(103, 17).Value = 3.25
$ws.Cells.Item(103, 19).Value = 1.825
$ws.Cells.Item(103, 20).Value = 1.975
$ws.Cells.Item(103, 21).Value = 2.5
$ws.Cells.Item(103, 22).Value = 1.95
$ws.Cells.Item(103, 23).Value = 1.75
$ws.Cells.Item(103, 24).Value = -1
$ws.Cells.Item(103, 25).Value = 2.3
$ws.Cells.Item(103, 27).Value = -0.5
$ws.Cells.Item(103, 28).Value = 0.4875
$ws.Cells.Item(103, 29).Value = -1
$ws.Cells.Item(103, 30).Value = 0.75

# Group: rows [114, 115]
$ws.Cells.Item(114, 2).Value = 6156940
$ws.Cells.Item(114, 5).Value = 'Universidad de Concepcion'
$ws.Cells.Item(114, 6).Value = 'Union San Felipe'
$ws.Cells.Item(114, 8).Value = 0
$ws.Cells.Item(114, 11).Value = 'D'
$ws.Cells.Item(114, 12).Value = 2.3
$ws.Cells.Item(114, 13).Value = 3.1
$ws.Cells.Item(114, 14).Value = 2.9
$ws.Cells.Item(114, 15).Value = 2.2
$ws.Cells.Item(114, 17).Value = 3.1
$ws.Cells.Item(114, 18).Value = -0.25
$ws.Cells.Item(114, 19).Value = 1.9
$ws.Cells.Item(114, 20).Value = 1.9
$ws.Cells.Item(114, 22).Value = 1.825
$ws.Cells.Item(114, 23).Value = 1.975
$ws.Cells.Item(114, 25).Value = 2.2
$ws.Cells.Item(114, 26).Value = -1
$ws.Cells.Item(114, 27).Value = -0.5
$ws.Cells.Item(114, 28).Value = 0.45
$ws.Cells.Item(114, 30).Value = 0.9750000000000001
$ws.Cells.Item(115, 2).Value = 6155961
$ws.Cells.Item(115, 5).Value = 'La Serena'
$ws.Cells.Item(115, 6).Value = 'Santiago Morning'
$ws.Cells.Item(115, 8).Value = 1
$ws.Cells.Item(115, 11).Value = 'A'
$ws.Cells.Item(115, 12).Value = 1.8
$ws.Cells.Item(115, 13).Value = 3.3
$ws.Cells.Item(115, 14).Value = 4
$ws.Cells.Item(115, 15).Value = 2.05
$ws.Cells.Item(115, 17).Value = 3.3
$ws.Cells.Item(115, 18).Value = -0.5
$ws.Cells.Item(115, 19).Value = 2.05
$ws.Cells.Item(115, 20).Value = 1.75
$ws.Cells.Item(115, 22).Value = 1.925
$ws.Cells.Item(115, 23).Value = 1.875
$ws.Cells.Item(115, 25).Value = -1
$ws.Cells.Item(115, 26).Value = 2.3
$ws.Cells.Item(115, 27).Value = -1
$ws.Cells.Item(115, 28).Value = 0.75
$ws.Cells.Item(115, 30).Value = 0.875

# Group: rows [116, 117]
$ws.Cells.Item(116, 2).Value = 6156942
$ws.Cells.Item(116, 5).Value = 'Deportes Recoleta'
$ws.Cells.Item(116, 6).Value = 'Rangers de Talca'
$ws.Cells.Item(116, 7).Value = 1
$ws.Cells.Item(116, 8).Value = 2
$ws.Cells.Item(116, 9).Value = 1
$ws.Cells.Item(116, 12).Value = 2.2
$ws.Cells.Item(116, 13).Value = 3.1
$ws.Cells.Item(116, 14).Value = 3
$ws.Cells.Item(116, 15).Value = 2.05
$ws.Cells.Item(116, 16).Value = 3.2
$ws.Cells.Item(116, 17).Value = 3.2
$ws.Cells.Item(116, 18).Value = -0.25
$ws.Cells.Item(116, 19).Value = 1.85
$ws.Cells.Item(116, 20).Value = 1.95
$ws.Cells.Item(116, 22).Value = 1.8
$ws.Cells.Item(116, 23).Value = 2
$ws.Cells.Item(116, 26).Value = 2.2
$ws.Cells.Item(116, 28).Value = 0.95
$ws.Cells.Item(116, 29).Value = 0.8
$ws.Cells.Item(117, 2).Value = 6155962
$ws.Cells.Item(117, 5).Value = 'Puerto Montt'
$ws.Cells.Item(117, 6).Value = 'Barnechea'
$ws.Cells.Item(117, 7).Value = 0
$ws.Cells.Item(117, 8).Value = 3
$ws.Cells.Item(117, 9).Value = 0
$ws.Cells.Item(117, 12).Value = 2.3
$ws.Cells.Item(117, 13).Value = 3.2
$ws.Cells.Item(117, 14).Value = 2.75
$ws.Cells.Item(117, 15).Value = 1.833
$ws.Cells.Item(117, 16).Value = 3.3
$ws.Cells.Item(117, 17).Value = 3.8
$ws.Cells.Item(117, 18).Value = -0.5
$ws.Cells.Item(117, 19).Value = 1.9
$ws.Cells.Item(117, 20).Value = 1.9
$ws.Cells.Item(117, 22).Value = 1.95
$ws.Cells.Item(117, 23).Value = 1.85
$ws.Cells.Item(117, 26).Value = 2.8
$ws.Cells.Item(117, 28).Value = 0.8999999999999999
$ws.Cells.Item(117, 29).Value = 0.95

# Group: rows [119, 120, 121]
$ws.Cells.Item(119, 2).Value = 6156939
$ws.Cells.Item(119, 5).Value = 'Club Deportes Santa Cruz'
$ws.Cells.Item(119, 6).Value = 'Deportes Iquique'
$ws.Cells.Item(119, 7).Value = 1
$ws.Cells.Item(119, 8).Value = 2
$ws.Cells.Item(119, 9).Value = 1
$ws.Cells.Item(119, 10).Value = 1
$ws.Cells.Item(119, 11).Value = 'A'
$ws.Cells.Item(119, 12).Value = 2.375
$ws.Cells.Item(119, 13).Value = 3.1
$ws.Cells.Item(119, 14).Value = 2.75
$ws.Cells.Item(119, 15).Value = 2.45
$ws.Cells.Item(119, 17).Value = 2.625
$ws.Cells.Item(119, 19).Value = 1.775
$ws.Cells.Item(119, 20).Value = 2.025
$ws.Cells.Item(119, 21).Value = 2.5
$ws.Cells.Item(119, 22).Value = 1.85
$ws.Cells.Item(119, 23).Value = 1.95
$ws.Cells.Item(119, 24).Value = -1
$ws.Cells.Item(119, 26).Value = 1.625
$ws.Cells.Item(119, 27).Value = -1
$ws.Cells.Item(119, 28).Value = 1.025
$ws.Cells.Item(119, 29).Value = 0.8500000000000001
$ws.Cells.Item(120, 2).Value = 6156941
$ws.Cells.Item(120, 5).Value = 'Cobreloa'
$ws.Cells.Item(120, 6).Value = 'San Luis Quillota'
$ws.Cells.Item(120, 8).Value = 0
$ws.Cells.Item(120, 9).Value = 0
$ws.Cells.Item(120, 10).Value = 0
$ws.Cells.Item(120, 11).Value = 'H'
$ws.Cells.Item(120, 12).Value = 1.8
$ws.Cells.Item(120, 13).Value = 3.3
$ws.Cells.Item(120, 14).Value = 4
$ws.Cells.Item(120, 15).Value = 1.615
$ws.Cells.Item(120, 16).Value = 3.6
$ws.Cells.Item(120, 17).Value = 4.75
$ws.Cells.Item(120, 18).Value = -0.75
$ws.Cells.Item(120, 19).Value = 1.875
$ws.Cells.Item(120, 20).Value = 1.925
$ws.Cells.Item(120, 21).Value = 2.25
$ws.Cells.Item(120, 22).Value = 1.825
$ws.Cells.Item(120, 23).Value = 1.975
$ws.Cells.Item(120, 24).Value = 0.615
$ws.Cells.Item(120, 26).Value = -1
$ws.Cells.Item(120, 27).Value = 0.4375
$ws.Cells.Item(120, 28).Value = -0.5
$ws.Cells.Item(120, 29).Value = -1
$ws.Cells.Item(120, 30).Value = 0.9750000000000001
$ws.Cells.Item(121, 2).Value = 6156184
$ws.Cells.Item(121, 5).Value = 'San Marcos De Arica'
$ws.Cells.Item(121, 6).Value = 'CD Antofagasta'
$ws.Cells.Item(121, 7).Value = 4
$ws.Cells.Item(121, 8).Value = 3
$ws.Cells.Item(121, 9).Value = 2
$ws.Cells.Item(121, 10).Value = 3
$ws.Cells.Item(121, 12).Value = 2.5
$ws.Cells.Item(121, 13).Value = 3.2
$ws.Cells.Item(121, 14).Value = 2.5
$ws.Cells.Item(121, 15).Value = 2.7
$ws.Cells.Item(121, 16).Value = 3.2
$ws.Cells.Item(121, 17).Value = 2.375
$ws.Cells.Item(121, 18).Value = 0
$ws.Cells.Item(121, 19).Value = 2.05
$ws.Cells.Item(121, 20).Value = 1.75
$ws.Cells.Item(121, 21).Value = 2.75
$ws.Cells.Item(121, 22).Value = 1.875
$ws.Cells.Item(121, 23).Value = 1.925
$ws.Cells.Item(121, 24).Value = 1.7
$ws.Cells.Item(121, 27).Value = 1.05
$ws.Cells.Item(121, 28).Value = -1
$ws.Cells.Item(121, 29).Value = 0.875
$ws.Cells.Item(121, 30).Value = -1

# Group: rows [122, 124]
$ws.Cells.Item(122, 2).Value = 7327838
$ws.Cells.Item(122, 5).Value = 'Barnechea'
$ws.Cells.Item(122, 6).Value = 'San Marcos De Arica'
$ws.Cells.Item(122, 7).Value = 3
$ws.Cells.Item(122, 8).Value = 3
$ws.Cells.Item(122, 9).Value = 1
$ws.Cells.Item(122, 10).Value = 2
$ws.Cells.Item(122, 11).Value = 'D'
$ws.Cells.Item(122, 12).Value = 2
$ws.Cells.Item(122, 13).Value = 3.3
$ws.Cells.Item(122, 14).Value = 3.2
$ws.Cells.Item(122, 15).Value = 2.1
$ws.Cells.Item(122, 16).Value = 3.3
$ws.Cells.Item(122, 17).Value = 3
$ws.Cells.Item(122, 18).Value = -0.25
$ws.Cells.Item(122, 19).Value = 1.9
$ws.Cells.Item(122, 20).Value = 1.9
$ws.Cells.Item(122, 21).Value = 3
$ws.Cells.Item(122, 22).Value = 2
$ws.Cells.Item(122, 23).Value = 1.8
$ws.Cells.Item(122, 25).Value = 2.3
$ws.Cells.Item(122, 26).Value = -1
$ws.Cells.Item(122, 27).Value = -0.5
$ws.Cells.Item(122, 28).Value = 0.45
$ws.Cells.Item(122, 29).Value = 1
$ws.Cells.Item(122, 30).Value = -1
$ws.Cells.Item(124, 2).Value = 7327856
$ws.Cells.Item(124, 5).Value = 'Union San Felipe'
$ws.Cells.Item(124, 6).Value = 'Puerto Montt'
$ws.Cells.Item(124, 7).Value = 0
$ws.Cells.Item(124, 8).Value = 1
$ws.Cells.Item(124, 9).Value = 0
$ws.Cells.Item(124, 10).Value = 0
$ws.Cells.Item(124, 11).Value = 'A'
$ws.Cells.Item(124, 12).Value = 1.727
$ws.Cells.Item(124, 13).Value = 3.5
$ws.Cells.Item(124, 14).Value = 4
$ws.Cells.Item(124, 15).Value = 1.8
$ws.Cells.Item(124, 16).Value = 3.4
$ws.Cells.Item(124, 17).Value = 3.75
$ws.Cells.Item(124, 18).Value = -0.5
$ws.Cells.Item(124, 19).Value = 1.85
$ws.Cells.Item(124, 20).Value = 1.95
$ws.Cells.Item(124, 21).Value = 2.25
$ws.Cells.Item(124, 22).Value = 1.75
$ws.Cells.Item(124, 23).Value = 1.95
$ws.Cells.Item(124, 25).Value = -1
$ws.Cells.Item(124, 26).Value = 2.75
$ws.Cells.Item(124, 27).Value = -1
$ws.Cells.Item(124, 28).Value = 0.95
$ws.Cells.Item(124, 29).Value = -1
$ws.Cells.Item(124, 30).Value = 0.95

# Group: rows [138, 139]
$ws.Cells.Item(138, 2).Value = 7503210
$ws.Cells.Item(138, 5).Value = 'Santiago Wanderers'
$ws.Cells.Item(138, 6).Value = 'Deportes Temuco'
$ws.Cells.Item(138, 7).Value = 1
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 13).Value = 3.2
$ws.Cells.Item(138, 14).Value = 3
$ws.Cells.Item(138, 16).Value = 3.2
$ws.Cells.Item(138, 17).Value = 3.1
$ws.Cells.Item(138, 19).Value = 1.925
$ws.Cells.Item(138, 20).Value = 1.875
$ws.Cells.Item(138, 21).Value = 2.5
$ws.Cells.Item(138, 22).Value = 1.975
$ws.Cells.Item(138, 23).Value = 1.825
$ws.Cells.Item(138, 27).Value = 0.925
$ws.Cells.Item(138, 29).Value = -1
$ws.Cells.Item(138, 30).Value = 0.825
$ws.Cells.Item(139, 2).Value = 7503209
$ws.Cells.Item(139, 5).Value = 'Deportes Iquique'
$ws.Cells.Item(139, 6).Value = 'CD Antofagasta'
$ws.Cells.Item(139, 7).Value = 3
$ws.Cells.Item(139, 9).Value = 1
$ws.Cells.Item(139, 13).Value = 3.6
$ws.Cells.Item(139, 14).Value = 2.75
$ws.Cells.Item(139, 16).Value = 3.6
$ws.Cells.Item(139, 17).Value = 2.875
$ws.Cells.Item(139, 19).Value = 1.975
$ws.Cells.Item(139, 20).Value = 1.825
$ws.Cells.Item(139, 21).Value = 2.75
$ws.Cells.Item(139, 22).Value = 1.85
$ws.Cells.Item(139, 23).Value = 1.95
$ws.Cells.Item(139, 27).Value = 0.9750000000000001
$ws.Cells.Item(139, 29).Value = 0.425
$ws.Cells.Item(139, 30).Value = -0.5

# Group: rows [177, 178]
$ws.Cells.Item(177, 2).Value = 7792884
$ws.Cells.Item(177, 5).Value = 'Santiago Wanderers'
$ws.Cells.Item(177, 6).Value = 'Deportes Temuco'
$ws.Cells.Item(177, 9).Value = 1
$ws.Cells.Item(177, 10).Value = 1
$ws.Cells.Item(177, 12).Value = 1.727
$ws.Cells.Item(177, 13).Value = 3.6
$ws.Cells.Item(177, 14).Value = 4
$ws.Cells.Item(177, 15).Value = 1.727
$ws.Cells.Item(177, 16).Value = 3.3
$ws.Cells.Item(177, 17).Value = 4.2
$ws.Cells.Item(177, 19).Value = 1.8
$ws.Cells.Item(177, 20).Value = 2
$ws.Cells.Item(177, 22).Value = 1.95
$ws.Cells.Item(177, 23).Value = 1.75
$ws.Cells.Item(177, 24).Value = 0.7270000000000001
$ws.Cells.Item(177, 27).Value = 0.8
$ws.Cells.Item(177, 29).Value = 0.95
$ws.Cells.Item(178, 2).Value = 7793501
$ws.Cells.Item(178, 5).Value = 'Magallanes'
$ws.Cells.Item(178, 6).Value = 'Santiago Morning'
$ws.Cells.Item(178, 9).Value = 2
$ws.Cells.Item(178, 10).Value = 0
$ws.Cells.Item(178, 12).Value = 2.3
$ws.Cells.Item(178, 13).Value = 3.4
$ws.Cells.Item(178, 14).Value = 2.625
$ws.Cells.Item(178, 15).Value = 1.909
$ws.Cells.Item(178, 16).Value = 3.6
$ws.Cells.Item(178, 17).Value = 3.3
$ws.Cells.Item(178, 19).Value = 1.925
$ws.Cells.Item(178, 20).Value = 1.875
$ws.Cells.Item(178, 22).Value = 1.9
$ws.Cells.Item(178, 23).Value = 1.9
$ws.Cells.Item(178, 24).Value = 0.909
$ws.Cells.Item(178, 27).Value = 0.925
$ws.Cells.Item(178, 29).Value = 0.8999999999999999

# Group: rows [194, 195]
$ws.Cells.Item(194, 2).Value = 7793509
$ws.Cells.Item(194, 5).Value = 'Barnechea'
$ws.Cells.Item(194, 6).Value = 'San Luis Quillota'
$ws.Cells.Item(194, 7).Value = 2
$ws.Cells.Item(194, 8).Value = 3
$ws.Cells.Item(194, 12).Value = 1.833
$ws.Cells.Item(194, 13).Value = 3.4
$ws.Cells.Item(194, 14).Value = 3.6
$ws.Cells.Item(194, 15).Value = 2.1
$ws.Cells.Item(194, 17).Value = 3.1
$ws.Cells.Item(194, 19).Value = 1.875
$ws.Cells.Item(194, 20).Value = 1.925
$ws.Cells.Item(194, 21).Value = 2.5
$ws.Cells.Item(194, 22).Value = 1.9
$ws.Cells.Item(194, 23).Value = 1.9
$ws.Cells.Item(194, 26).Value = 2.1
$ws.Cells.Item(194, 28).Value = 0.925
$ws.Cells.Item(194, 29).Value = 0.8999999999999999
$ws.Cells.Item(194, 30).Value = -1
$ws.Cells.Item(195, 2).Value = 7793510
$ws.Cells.Item(195, 5).Value = 'Deportes Temuco'
$ws.Cells.Item(195, 6).Value = 'Santiago Morning'
$ws.Cells.Item(195, 7).Value = 0
$ws.Cells.Item(195, 8).Value = 1
$ws.Cells.Item(195, 12).Value = 2.05
$ws.Cells.Item(195, 13).Value = 3.25
$ws.Cells.Item(195, 14).Value = 3.1
$ws.Cells.Item(195, 15).Value = 2.15
$ws.Cells.Item(195, 17).Value = 2.9
$ws.Cells.Item(195, 19).Value = 1.975
$ws.Cells.Item(195, 20).Value = 1.825
$ws.Cells.Item(195, 21).Value = 2.25
$ws.Cells.Item(195, 22).Value = 1.75
$ws.Cells.Item(195, 23).Value = 1.95
$ws.Cells.Item(195, 26).Value = 1.9
$ws.Cells.Item(195, 28).Value = 0.825
$ws.Cells.Item(195, 29).Value = -1
$ws.Cells.Item(195, 30).Value = 0.95

# Group: rows [201, 202]
$ws.Cells.Item(201, 2).Value = 7793513
$ws.Cells.Item(201, 5).Value = 'San Luis Quillota'
$ws.Cells.Item(201, 6).Value = 'La Serena'
$ws.Cells.Item(201, 7).Value = 0
$ws.Cells.Item(201, 8).Value = 1
$ws.Cells.Item(201, 9).Value = 0
$ws.Cells.Item(201, 10).Value = 0
$ws.Cells.Item(201, 12).Value = 2.75
$ws.Cells.Item(201, 13).Value = 3.4
$ws.Cells.Item(201, 14).Value = 2.2
$ws.Cells.Item(201, 15).Value = 1.8
$ws.Cells.Item(201, 16).Value = 3.6
$ws.Cells.Item(201, 17).Value = 3.5
$ws.Cells.Item(201, 18).Value = -0.5
$ws.Cells.Item(201, 19).Value = 1.85
$ws.Cells.Item(201, 20).Value = 1.95
$ws.Cells.Item(201, 21).Value = 2.25
$ws.Cells.Item(201, 22).Value = 1.85
$ws.Cells.Item(201, 23).Value = 1.95
$ws.Cells.Item(201, 26).Value = 2.5
$ws.Cells.Item(201, 28).Value = 0.95
$ws.Cells.Item(201, 29).Value = -1
$ws.Cells.Item(201, 30).Value = 0.95
$ws.Cells.Item(202, 2).Value = 7792892
$ws.Cells.Item(202, 5).Value = 'Union San Felipe'
$ws.Cells.Item(202, 6).Value = 'Deportes Temuco'
$ws.Cells.Item(202, 7).Value = 2
$ws.Cells.Item(202, 8).Value = 3
$ws.Cells.Item(202, 9).Value = 1
$ws.Cells.Item(202, 10).Value = 2
$ws.Cells.Item(202, 12).Value = 2.375
$ws.Cells.Item(202, 13).Value = 3.3
$ws.Cells.Item(202, 14).Value = 2.625
$ws.Cells.Item(202, 15).Value = 2.1
$ws.Cells.Item(202, 16).Value = 3.2
$ws.Cells.Item(202, 17).Value = 3.1
$ws.Cells.Item(202, 18).Value = -0.25
$ws.Cells.Item(202, 19).Value = 1.875
$ws.Cells.Item(202, 20).Value = 1.925
$ws.Cells.Item(202, 21).Value = 2.5
$ws.Cells.Item(202, 22).Value = 1.95
$ws.Cells.Item(202, 23).Value = 1.85
$ws.Cells.Item(202, 26).Value = 2.1
$ws.Cells.Item(202, 28).Value = 0.925
$ws.Cells.Item(202, 29).Value = 0.95
$ws.Cells.Item(202, 30).Value = -1

